$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$shp = $s.Shapes.Item(2)
$tr = $shp.TextFrame.TextRange

$word = $tr.Characters(1, 6)
$word.Text = "Class"
$word.LanguageID = "en-US"

$space = $tr.Characters(6, 1)
$space.Text = " "
